# "Generate Report for Handback"
#
# The handback report previously contained two file entries
# (4f3491b5-...  and a53b297b-...). This run only produced a handback
# report for 4f3491b5-..., so the a53b297b-... row is removed from every
# sheet (Overview, zh-cn, de-de), and the remaining row's handoff/handback
# timestamps are refreshed to the new run's values.

$wb = $excel.ActiveWorkbook

$HYPER_COLOR = 15570276   # BGR for RGB(100,149,237) / #6495ED -> matches the
                           # workbook's existing "HyperLink" font color

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Delete()
$wsOverview.Rows.Item(3).Delete()

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/785d70189b44bc0ae84287e30897497a8b506391/e2e/4f3491b5-b04d-497b-a309-93df6e5fb3f9.md",
    [Type]::Missing,
    [Type]::Missing,
    "4f3491b5-b04d-497b-a309-93df6e5fb3f9.md")
$wsOverview.Range("A2").Font.Underline = $true
$wsOverview.Range("A2").Font.Color = $HYPER_COLOR

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Rows.Item(3).Delete()

# Refresh the Correspond Handoff/Handback Datetime values for the
# remaining row to this run's timestamps.
$wsZhCn.Range("E2").Value = "2016-03-23 20:56:24"
$wsZhCn.Range("H2").Value = "2016-03-23 20:56:48"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/785d70189b44bc0ae84287e30897497a8b506391/e2e/4f3491b5-b04d-497b-a309-93df6e5fb3f9.md",
    [Type]::Missing,
    [Type]::Missing,
    "4f3491b5-b04d-497b-a309-93df6e5fb3f9.md")
$wsZhCn.Range("A2").Font.Underline = $true
$wsZhCn.Range("A2").Font.Color = $HYPER_COLOR

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("D2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f6fb328f726ddb727f86da016d3006d5b35d352b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/4f3491b5-b04d-497b-a309-93df6e5fb3f9.548891f56019860c5abafc3df1f7916c47c372a6.zh-cn.xlf",
    [Type]::Missing,
    [Type]::Missing,
    "4f3491b5-b04d-497b-a309-93df6e5fb3f9.548891f56019860c5abafc3df1f7916c47c372a6.zh-cn.xlf")
$wsZhCn.Range("D2").Font.Underline = $true
$wsZhCn.Range("D2").Font.Color = $HYPER_COLOR

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("F2"),
    "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/ddf849e67ce333f8b559ca4e87b73a4cfb7449cb/e2e/4f3491b5-b04d-497b-a309-93df6e5fb3f9.md",
    [Type]::Missing,
    [Type]::Missing,
    "4f3491b5-b04d-497b-a309-93df6e5fb3f9.md")
$wsZhCn.Range("F2").Font.Underline = $true
$wsZhCn.Range("F2").Font.Color = $HYPER_COLOR

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/de4efd86f4bc1a86bc21041a461e611dc43bd40d/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/4f3491b5-b04d-497b-a309-93df6e5fb3f9.548891f56019860c5abafc3df1f7916c47c372a6.zh-cn.xlf",
    [Type]::Missing,
    [Type]::Missing,
    "4f3491b5-b04d-497b-a309-93df6e5fb3f9.548891f56019860c5abafc3df1f7916c47c372a6.zh-cn.xlf")
$wsZhCn.Range("G2").Font.Underline = $true
$wsZhCn.Range("G2").Font.Color = $HYPER_COLOR

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Rows.Item(3).Delete()

$wsDeDe.Range("E2").Value = "2016-03-23 20:56:29"
$wsDeDe.Range("H2").Value = "2016-03-23 20:56:54"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/785d70189b44bc0ae84287e30897497a8b506391/e2e/4f3491b5-b04d-497b-a309-93df6e5fb3f9.md",
    [Type]::Missing,
    [Type]::Missing,
    "4f3491b5-b04d-497b-a309-93df6e5fb3f9.md")
$wsDeDe.Range("A2").Font.Underline = $true
$wsDeDe.Range("A2").Font.Color = $HYPER_COLOR

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("D2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/10c75d53c48a71975083f97e1d7f68a66b707e8d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/4f3491b5-b04d-497b-a309-93df6e5fb3f9.548891f56019860c5abafc3df1f7916c47c372a6.de-de.xlf",
    [Type]::Missing,
    [Type]::Missing,
    "4f3491b5-b04d-497b-a309-93df6e5fb3f9.548891f56019860c5abafc3df1f7916c47c372a6.de-de.xlf")
$wsDeDe.Range("D2").Font.Underline = $true
$wsDeDe.Range("D2").Font.Color = $HYPER_COLOR

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("F2"),
    "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/c42b002d72195fb30d42ac765e4188e826a7b0d3/e2e/4f3491b5-b04d-497b-a309-93df6e5fb3f9.md",
    [Type]::Missing,
    [Type]::Missing,
    "4f3491b5-b04d-497b-a309-93df6e5fb3f9.md")
$wsDeDe.Range("F2").Font.Underline = $true
$wsDeDe.Range("F2").Font.Color = $HYPER_COLOR

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/434a095bc8683c52cd639890c1ee8e264586742e/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/4f3491b5-b04d-497b-a309-93df6e5fb3f9.548891f56019860c5abafc3df1f7916c47c372a6.de-de.xlf",
    [Type]::Missing,
    [Type]::Missing,
    "4f3491b5-b04d-497b-a309-93df6e5fb3f9.548891f56019860c5abafc3df1f7916c47c372a6.de-de.xlf")
$wsDeDe.Range("G2").Font.Underline = $true
$wsDeDe.Range("G2").Font.Color = $HYPER_COLOR
